$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (Df, SumsOfSqs, MeanSqs, F.Model, R2, Pr(>F))
$data = @(
    @(1, 1.99007531613154, 1.99007531613154, 2.79809603102731, 0.0122057646948616, 0.1038),
    @(1, 2.57658610808061, 2.57658610808061, 3.6227449806447, 0.0158030218737724, 0.0578),
    @(3, 24.8358273201153, 8.27860910670511, 11.6399329694349, 0.152326025962078, 0.0001),
    @(1, 0.145274435574549, 0.145274435574549, 0.204259516358933, 0.00089101430605583, 0.6475),
    @(3, 0.95191573660156, 0.317305245533853, 0.446139169183954, 0.00583840189168323, 0.7173),
    @(3, 38.4866673488115, 12.8288891162705, 18.0377412994118, 0.236050968393578, 0.0001),
    @(3, 0.175865615319353, 0.0586218717731176, 0.0824238285908722, 0.00107863973845861, 0.9698),
    @(132, 93.8816748304804, 0.711224809321821, $null, 0.575806163139512, $null),
    @(147, 163.043886711115, $null, $null, 1, $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $val = $rowData[$c]
        if ($null -eq $val) {
            $ws.Cells.Item($row, $c + 1).Value = $null
        } else {
            $ws.Cells.Item($row, $c + 1).Value = $val
        }
    }
}
